# B1--and-B2-PowerPoint.pptx edit
#
# 1. Slide 5's table switches from the custom "Table_0" style
#    ({5C0DC621-3EA0-4412-AF89-2DA6A7EA206F}) to the built-in table
#    style {61C430A0-938C-4D03-B4FD-118DBA27FC07}.
# 2. The deck's colour theme (the one actually driving every slide's
#    look, via the slide master) switches from the colourful "Integral"
#    / "Red Violet" palette to the plain default "Office Theme" palette
#    -- i.e. the 12 theme colours are swapped for the Office defaults.

function HexToComRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 5 -------------------------------------------
$slide5 = $p.Slides.Item(5)
$tableShape = $slide5.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{61C430A0-938C-4D03-B4FD-118DBA27FC07}", $true)

# --- 2. Swap the presentation's theme colours to the plain Office theme --
$officeThemeColors = @{
    1  = "000000" # dk1
    2  = "FFFFFF" # lt1
    3  = "44546A" # dk2
    4  = "E7E6E6" # lt2
    5  = "5B9BD5" # accent1
    6  = "ED7D31" # accent2
    7  = "A5A5A5" # accent3
    8  = "FFC000" # accent4
    9  = "4472C4" # accent5
    10 = "70AD47" # accent6
    11 = "0563C1" # hlink
    12 = "954F72" # folHlink
}

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = HexToComRgb($officeThemeColors[$i])
}
